$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37 and 38 swap content (PolygonEcosystemToken <-> ImmutableX) plus updated price/volume values
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.48%  "

$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "'0.435"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +15.45%  "

$ws.Range("D2").Value = "'57.817.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "'2.350.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("D5").Value = "'548.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "'132.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.10%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.568"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("E9").Value = "  +3.33%  "
$ws.Range("D10").Value = "'5.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.32%  "
$ws.Range("E11").Value = "  -1.19%  "
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("D13").Value = "'23.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").Value = "'2.768.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "'57.743.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "'0.0000136"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.63%  "
$ws.Range("D17").Value = "'2.360.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("E18").Value = "  +2.94%  "
$ws.Range("D19").Value = "'4.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("D20").Value = "'329.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'6.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.04%  "
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").Value = "'63.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.69%  "
$ws.Range("D24").Value = "'0.168"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").Value = "'8.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.33%  "
$ws.Range("E27").Value = "  -6.07%  "
$ws.Range("E28").Value = "  -1.07%  "
$ws.Range("D29").Value = "'170.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").Value = "'0.0₃0735"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").Value = "'6.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("E33").Value = "  -2.69%  "
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("D39").Value = "'40.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.11%  "
$ws.Range("E40").Value = "  -1.63%  "
$ws.Range("D41").Value = "'141.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.37%  "
$ws.Range("D42").Value = "'3.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").Value = "'287.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("D44").Value = "'0.423"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.86%  "
$ws.Range("D45").Value = "'0.0953"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("D46").Value = "'0.0512"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("D47").Value = "'0.567"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("D48").Value = "'18.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.85%  "
$ws.Range("D49").Value = "'0.0221"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.42%  "
$ws.Range("D50").Value = "'11.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("E51").Value = "  -0.07%  "
